# [BI-2053] - updated test files
#
# The "unit1" value in column O (rows 2-5) of the Template sheet is no
# longer needed, so clear those cells. Removing the last usage of the
# "unit1" shared string also drops it from the shared strings table and
# shifts the remaining string indices down automatically on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Template")
$ws.Activate()

$ws.Range("O2:O5").ClearContents()

# Reflect the author's updated view/selection state on the sheet.
$ws.Range("D1").Select()
$ws.Range("O2:O5").Select()
